# Update "Pais" (countries) worksheet:
#  - Costa Rica's stats improved enough that it now outranks Eslovaquia,
#    Principado de Andorra, Uruguay, Jordania, San Marino and Kuwait, so its
#    row moves up from row 80 to row 74; those six countries each shift down
#    one row (their own totals are unchanged).
#  - A handful of other countries got refreshed totals (Reino Unido, Austria,
#    Canada, Rumania, Argentina, Montenegro).
#  - The "last updated" timestamp banner in A1 moves from 19:42 to 20:12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 20:12"

# --- Independent total/stat refreshes ----------------------------------
# Reino Unido (row 12): Casos criticos
$ws.Range("G12").Value = 115

# Austria (row 15): Casos totales, Nuevos casos, Recuperados
$ws.Range("B15").Value = 6847
$ws.Range("C15").Value = 1259
$ws.Range("E15").Value = 6686

# Canada (row 17): Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes
$ws.Range("B17").Value = 3897
$ws.Range("C17").Value = 488
$ws.Range("D17").Value = 199
$ws.Range("E17").Value = 3661
$ws.Range("F17").Value = 120
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 37

# Rumania (row 36): Recuperados, Muertes hoy, Muertes
$ws.Range("E36").Value = 913
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 22

# Argentina (row 53): Casos activos, Recuperados
$ws.Range("D53").Value = 63
$ws.Range("E53").Value = 431

# Montenegro (row 110): Casos totales, Nuevos casos, Recuperados
$ws.Range("B110").Value = 69
$ws.Range("C110").Value = 16
$ws.Range("E110").Value = 68

# --- Costa Rica re-ranking (rows 74-80) --------------------------------
# Eslovaquia, Principado de Andorra, Uruguay, Jordania, San Marino and Kuwait
# each keep their existing totals but shift down one row (74->75, 75->76,
# ..., 79->80) to make room for Costa Rica's improved totals at row 74.
# Capture their current (pre-shift) name + totals first so the writes below
# can't clobber a value before it has been read.
$rows = @(74, 75, 76, 77, 78, 79)
$cols = @("A", "B", "C", "D", "E", "F", "G", "H")

$captured = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value()
    }
    $captured[$r] = $rowData
}

foreach ($r in $rows) {
    $destRow = $r + 1
    $rowData = $captured[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$col]
    }
}

# Put Costa Rica's updated totals into row 74 (its new ranked position)
$ws.Range("A74").Value = "Costa Rica"
$ws.Range("B74").Value = 231
$ws.Range("C74").Value = 30
$ws.Range("D74").Value = 2
$ws.Range("E74").Value = 227
$ws.Range("F74").Value = 4
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 2
